# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# de-de handback has completed (handoff/handback round-trip finished),
# while zh-cn additionally records its handback datetime.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: update Status-like cells (E2/F2/E3/F3) and widen columns
# E/F to match the new, longer status text.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------------
# Helper values shared between the zh-cn and de-de sheets.
# ---------------------------------------------------------------------------
$mdFile1 = "048510b4-dbe2-4a4a-8231-6b20506ae6b2.md"
$mdFile2 = "c0635734-0671-4fae-97fd-ad4939ff1ade.md"
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feaca489a27bfd622516297a8e2614bbcef6b607/e2e/048510b4-dbe2-4a4a-8231-6b20506ae6b2.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feaca489a27bfd622516297a8e2614bbcef6b607/e2e/c0635734-0671-4fae-97fd-ad4939ff1ade.md"

$hyperFontColor = 15570276  # OLE (BGR) form of RGB 6495ED, the workbook's HyperLink font color

function Apply-HyperlinkFont($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperFontColor
}

# ---------------------------------------------------------------------------
# zh-cn sheet: fill in the newly generated handoff target / handback file
# columns (I/J) for both rows, record the handback datetime, widen columns,
# and re-lay the hyperlinks so that A2, I2, A3, I3 end up with freshly
# assigned relationship ids in that order.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsZhCn.Range("I2").Value = $mdFile1
$wsZhCn.Range("J2").Value = "048510b4-dbe2-4a4a-8231-6b20506ae6b2.8aa429726f4b0885d68a7b3a74deae607013df4c.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-25 08:51:21"

$wsZhCn.Range("I3").Value = $mdFile2
$wsZhCn.Range("J3").Value = "c0635734-0671-4fae-97fd-ad4939ff1ade.f925a3374debe69e4f81f9109114979ba048505b.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-25 08:51:21"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $mdFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $mdFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $mdFile2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $mdFile2)

Apply-HyperlinkFont $wsZhCn.Range("A2")
Apply-HyperlinkFont $wsZhCn.Range("A3")
Apply-HyperlinkFont $wsZhCn.Range("I2")
Apply-HyperlinkFont $wsZhCn.Range("I3")

# ---------------------------------------------------------------------------
# de-de sheet: same treatment, but both rows share a single, later handback
# datetime.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

$wsDeDe.Range("I2").Value = $mdFile1
$wsDeDe.Range("J2").Value = "048510b4-dbe2-4a4a-8231-6b20506ae6b2.8aa429726f4b0885d68a7b3a74deae607013df4c.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-25 08:51:28"

$wsDeDe.Range("I3").Value = $mdFile2
$wsDeDe.Range("J3").Value = "c0635734-0671-4fae-97fd-ad4939ff1ade.f925a3374debe69e4f81f9109114979ba048505b.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-25 08:51:28"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $mdFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $mdFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $mdFile2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $mdFile2)

Apply-HyperlinkFont $wsDeDe.Range("A2")
Apply-HyperlinkFont $wsDeDe.Range("A3")
Apply-HyperlinkFont $wsDeDe.Range("I2")
Apply-HyperlinkFont $wsDeDe.Range("I3")
